$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged) - update values
$ws.Range("B3").Value = 0.9945325729021466
$ws.Range("C3").Value = 0.9945276457813638
$ws.Range("D3").Value = 0.988963261755489

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9881624143358607
$ws.Range("C4").Value = 0.9889261686743417
$ws.Range("D4").Value = 0.9712075185729506

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9972983235499356
$ws.Range("C5").Value = 0.997250406505402
$ws.Range("D5").Value = 0.9965031487786785
